$d = $word.ActiveDocument

# Pull the full package as flat OOXML so we can perform precise, surgical
# text-level edits that mirror the authored diff (the real edits live deep
# inside nested group-shape text boxes that the high level Shapes/GroupItems
# object model in this doc cannot reach).
$xml = $d.Content.WordOpenXML

# ---------------------------------------------------------------------
# 1) Bookmark: close the "_GoBack" bookmark immediately after it starts
#    instead of at the end of the paragraph.
# ---------------------------------------------------------------------
$xml = $xml.Replace(
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:r>',
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r>'
)
$xml = $xml.Replace(
    '</mc:Fallback></mc:AlternateContent></w:r><w:bookmarkEnd w:id="0"/></w:p>',
    '</mc:Fallback></mc:AlternateContent></w:r></w:p>'
)

# ---------------------------------------------------------------------
# 2) Two of the five identical "grey rectangle" fills (the callout text
#    box backgrounds "Text Box 8" / "Text Box 38" in the DrawingML
#    wordprocessingGroup) become transparent (noFill).
# ---------------------------------------------------------------------
$oldFill = '<a:solidFill><a:schemeClr val="bg1"><a:lumMod val="65000"/></a:schemeClr></a:solidFill><a:ln w="6350"><a:noFill/></a:ln>'
$newFill = '<a:noFill/><a:ln w="6350"><a:noFill/></a:ln>'
$xml = $xml.Replace($oldFill, $newFill)

# ---------------------------------------------------------------------
# 3) Split the "bit o' charge:  " run into two runs: a plain-font label
#    run ("bit o' charge:") and a Times-New-Roman run holding just the
#    trailing two spaces. Occurs 4x (DrawingML x2 + VML fallback x2).
# ---------------------------------------------------------------------
$oldRun = '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">bit o’ charge:  </w:t></w:r>'
$newRun = '<w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia" w:cstheme="minorHAnsi"/></w:rPr><w:t>bit o’ charge:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">  </w:t></w:r>'
$xml = $xml.Replace($oldRun, $newRun)

$d.Content.WordOpenXML = $xml

Write-Output "done"
